$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Adam9", "Itga3", "ECs", 3, 1, 11.915632, 35.746896, 0.2203762099850903, 0.2203762099850904, 3, 1, 8.269168666666666, 24.807506, 0.671680253471746, 0.671680253471746, 98.53237077793065, 886.7913370013758, 0.1480223485819282, 0.1480223485819282),
    @("ECs", "Adam9", "Itga3", "FAPs", 3, 1, 11.915632, 35.746896, 0.2203762099850903, 0.2203762099850904, 2, 0.6666666666666666, 0.121294, 0.363882, 0.009852354928133683, 0.009852354928133683, 1.445294667808, 13.007652010272, 0.002171224638490028, 0.002171224638490028),
    @("ECs", "Adam9", "Itga3", "sCs", 3, 1, 11.915632, 35.746896, 0.2203762099850903, 0.2203762099850904, 3, 1, 3.920705666666667, 11.762117, 0.3184673916001203, 0.3184673916001203, 46.71768590431467, 420.459173138832, 0.0701826367646721, 0.07018263676467211),
    @("FAPs", "Adam9", "Itga3", "ECs", 3, 1, 29.800487, 89.40146100000001, 0.5511514941691683, 0.5511514941691684, 3, 1, 8.269168666666666, 24.807506, 0.671680253471746, 0.671680253471746, 246.4252533518073, 2217.827280166266, 0.3701975753048785, 0.3701975753048786),
    @("FAPs", "Adam9", "Itga3", "FAPs", 3, 1, 29.800487, 89.40146100000001, 0.5511514941691683, 0.5511514941691684, 2, 0.6666666666666666, 0.121294, 0.363882, 0.009852354928133683, 0.009852354928133683, 3.614620270178, 32.531582431602, 0.005430140139725849, 0.00543014013972585),
    @("FAPs", "Adam9", "Itga3", "sCs", 3, 1, 29.800487, 89.40146100000001, 0.5511514941691683, 0.5511514941691684, 3, 1, 3.920705666666667, 11.762117, 0.3184673916001203, 0.3184673916001203, 116.8389382503264, 1051.550444252937, 0.1755237787245639, 0.175523778724564),
    @("sCs", "Adam9", "Itga3", "ECs", 3, 1, 12.35338333333333, 37.06015, 0.2284722958457413, 0.2284722958457413, 3, 1, 8.269168666666666, 24.807506, 0.671680253471746, 0.671680253471746, 102.1522103873222, 919.3698934858999, 0.1534603295849392, 0.1534603295849392),
    @("sCs", "Adam9", "Itga3", "FAPs", 3, 1, 12.35338333333333, 37.06015, 0.2284722958457413, 0.2284722958457413, 2, 0.6666666666666666, 0.121294, 0.363882, 0.009852354928133683, 0.009852354928133683, 1.498391278033333, 13.4855215023, 0.002250990149917806, 0.002250990149917806),
    @("sCs", "Adam9", "Itga3", "sCs", 3, 1, 12.35338333333333, 37.06015, 0.2284722958457413, 0.2284722958457413, 3, 1, 3.920705666666667, 11.762117, 0.3184673916001203, 0.3184673916001203, 48.43398003750556, 435.90582033755, 0.07276097611088422, 0.07276097611088422),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
